# Apply text corrections (header rename + Spanish preposition capitalization)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textChanges = @(
    @('A1', 'mx_state'),
    @('B1', 'mx_municipality'),
    @('C1', 'n_matriculas'),
    @('D1', 'pct_matriculas'),
    @('B6', 'Pabellón De Arteaga'),
    @('B7', 'Rincón De Romos'),
    @('B8', 'San Francisco De Los Romo'),
    @('B28', 'Comitán De Domínguez'),
    @('B59', 'Hidalgo Del Parral'),
    @('B67', 'Valle De Zaragoza'),
    @('B80', 'Villa De Álvarez'),
    @('A82', 'Ciudad De México'),
    @('B86', 'Cuajimalpa De Morelos'),
    @('B99', 'Coneto De Comonfort'),
    @('B116', 'San Juan Del Río'),
    @('A124', 'Estado De México'),
    @('B124', 'Acambay De Ruíz Castañeda'),
    @('B125', 'Almoloya De Juárez'),
    @('B129', 'Atizapán De Zaragoza'),
    @('B133', 'Chapa De Mota'),
    @('B136', 'Coacalco De Berriozábal'),
    @('B138', 'Ecatepec De Morelos'),
    @('B148', 'San Felipe Del Progreso'),
    @('B153', 'Tenango Del Valle'),
    @('B158', 'Tlalnepantla De Baz'),
    @('B167', 'San Miguel De Allende'),
    @('B168', 'Apaseo El Alto'),
    @('B174', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B178', 'Jaral Del Progreso'),
    @('B185', 'Purísima Del Rincón'),
    @('B189', 'San Diego De La Unión'),
    @('B191', 'San Francisco Del Rincón'),
    @('B193', 'San Luis De La Paz'),
    @('B194', 'Santa Cruz De Juventino Rosas'),
    @('B195', 'Silao De La Victoria'),
    @('B198', 'Valle De Santiago'),
    @('B202', 'Acapulco De Juárez'),
    @('B203', 'Ajuchitlán Del Progreso'),
    @('B206', 'Atoyac De Álvarez'),
    @('B207', 'Ayutla De Los Libres'),
    @('B210', 'Buenavista De Cuéllar'),
    @('B211', 'Chilapa De Álvarez'),
    @('B212', 'Chilpancingo De Los Bravo'),
    @('B213', 'Coahuayutla De José María Izazaga'),
    @('B216', 'Coyuca De Benítez'),
    @('B219', 'Cuetzala Del Progreso'),
    @('B220', 'Cutzamala De Pinzón'),
    @('B224', 'Iguala De La Independencia'),
    @('B226', 'Ixcateopan De Cuauhtémoc'),
    @('B229', 'Mártir De Cuilapan'),
    @('B238', 'Taxco De Alarcón'),
    @('B239', 'Técpan De Galeana'),
    @('B243', 'Tlalixtaquilla De Maldonado'),
    @('B250', 'Agua Blanca De Iturbide'),
    @('B251', 'Atotonilco El Grande'),
    @('B252', 'Cuautepec De Hinojosa'),
    @('B259', 'Omitlán De Juárez'),
    @('B260', 'Pachuca De Soto'),
    @('B265', 'Tenango De Doria'),
    @('B267', 'Tepehuacán De Guerrero'),
    @('B268', 'Tulancingo De Bravo'),
    @('B276', 'Atotonilco El Alto'),
    @('B278', 'Autlán De Navarro'),
    @('B293', 'Huejuquilla El Alto'),
    @('B294', 'Ixtlahuacán Del Río'),
    @('B296', 'Jilotlán De Los Dolores'),
    @('B298', 'Lagos De Moreno'),
    @('B304', 'Ojuelos De Jalisco'),
    @('B307', 'San Cristóbal De La Barranca'),
    @('B310', 'San Miguel El Alto'),
    @('B311', 'Santa María De Los Ángeles'),
    @('B313', 'Tamazula De Gordiano'),
    @('B317', 'Teocuitatlán De Corona'),
    @('B318', 'Tepatitlán De Morelos'),
    @('B320', 'Tizapán El Alto'),
    @('B324', 'Unión De San Antonio'),
    @('B325', 'Valle De Juárez'),
    @('B328', 'Zapotlán El Grande'),
    @('B396', 'Puente De Ixtla'),
    @('B400', 'Tlaltizapán De Zapata'),
    @('B405', 'Amatlán De Cañas'),
    @('B409', 'Santa María Del Oro'),
    @('B421', 'San Nicolás De Los Garza'),
    @('B424', 'Acatlán De Pérez Figueroa'),
    @('B430', 'Guadalupe De Ramírez'),
    @('B431', 'Guevea De Humboldt'),
    @('B432', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B433', 'Heroica Ciudad De Huajuapan De León'),
    @('B434', 'Ixtlán De Juárez'),
    @('B435', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B440', 'Miahuatlán De Porfirio Díaz'),
    @('B442', 'Oaxaca De Juárez'),
    @('B447', 'San Antonino El Alto'),
    @('B449', 'San Francisco Del Mar'),
    @('B473', 'Santa Lucía Del Camino'),
    @('B474', 'Santa María Del Tule'),
    @('B490', 'Villa De Tamazulápam Del Progreso'),
    @('B491', 'Villa De Zaachila'),
    @('B492', 'Villa Sola De Vega'),
    @('B493', 'Zimatlán De Álvarez'),
    @('B509', 'Los Reyes De Juárez'),
    @('B518', 'Tepanco De López'),
    @('B519', 'Tlacotepec De Benito Juárez'),
    @('B529', 'Amealco De Bonfil'),
    @('B530', 'Cadereyta De Montes'),
    @('B532', 'Landa De Matamoros'),
    @('B535', 'Pinal De Amoles'),
    @('B537', 'San Juan Del Río'),
    @('B555', 'Santa María Del Río'),
    @('B560', 'Villa De Guadalupe'),
    @('B561', 'Villa De Ramos'),
    @('B601', 'Soto La Marina'),
    @('B610', 'Tetla De La Solidaridad'),
    @('B614', 'Alto Lucero De Gutiérrez Barrios'),
    @('B615', 'Amatlán De Los Reyes'),
    @('B623', 'Cosamaloapan De Carpio'),
    @('B636', 'Juchique De Ferrer'),
    @('B649', 'Paso De Ovejas'),
    @('B650', 'Paso Del Macho'),
    @('B653', 'Poza Rica De Hidalgo'),
    @('B656', 'Soledad De Doblado'),
    @('B676', 'El Plateado De Joaquín Amaro'),
    @('B689', 'Moyahua De Estrada'),
    @('B690', 'Nochistlán De Mejía'),
    @('B698', 'Tlaltenango De Sánchez Román'),
    @('B700', 'Villa De Cos'),
)

foreach ($change in $textChanges) {
    $ws.Range($change[0]).Value = $change[1]
}

# Fix floating point rounding on percentage column (D20, D76, D505)
$ws.Range("D20").Value = 0.009393346379647748
$ws.Range("D76").Value = 0.009393346379647748
$ws.Range("D505").Value = 0.009393346379647748

# Remove trailing footer/metadata rows (708:712) no longer part of the clean dataset
$ws.Rows("708:712").Delete()

Write-Output "edit complete"
